$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 8-15 and add new rows 16-17.
# Columns: A=index, B=name, C=from_bus, D=to_bus, E=in_service

$data = @(
    @{ row = 8;  a = 6;  name = "line7";  c = 14; d = 11; e = $true  },
    @{ row = 9;  a = 7;  name = "line8";  c = 16; d = 9;  e = $true  },
    @{ row = 10; a = 8;  name = "extr1";  c = 5;  d = 12; e = $true  },
    @{ row = 11; a = 9;  name = "extr2";  c = 5;  d = 9;  e = $true  },
    @{ row = 12; a = 10; name = "extr3";  c = 10; d = 11; e = $true  },
    @{ row = 13; a = 11; name = "extr4";  c = 7;  d = 8;  e = $false },
    @{ row = 14; a = 12; name = "extr5";  c = 9;  d = 11; e = $false },
    @{ row = 15; a = 13; name = "extr6";  c = 7;  d = 11; e = $true  },
    @{ row = 16; a = 14; name = "extr7";  c = 5;  d = 7;  e = $true  },
    @{ row = 17; a = 15; name = "extr8";  c = 8;  d = 5;  e = $true  }
)

# Rows 16-17 are brand new; clone the formatting of row 15 (the last existing
# data row) into them before writing values, so they match the look of the
# rest of the table (bold/centered/bordered index column, etc.).
$ws.Range("A15:E15").Copy() | Out-Null
$ws.Range("A16:E17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.a
    $ws.Cells.Item($r, 2).Value = $item.name
    $ws.Cells.Item($r, 3).Value = $item.c
    $ws.Cells.Item($r, 4).Value = $item.d
    $ws.Cells.Item($r, 5).Value = $item.e
}
